$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style from an existing header cell (e.g. E1) to the new headers
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) # xlPasteFormats

# Boolean values for rows 2-25, columns F, G, H - all FALSE by default
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}

# H14 is TRUE per the diff
$ws.Cells.Item(14, 8).Value = $true
